$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (MuSCs -> Fgf8 -> Fgfr2 -> ECs) with new TPM-derived values ---
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Fgf8"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07200033333333333
$ws.Range("H2").Value = 0.216001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.106124
$ws.Range("N2").Value = 0.318372
$ws.Range("O2").Value = 0.08094716512538251
$ws.Range("P2").Value = 0.08094716512538253
$ws.Range("Q2").Value = 0.007640963374666666
$ws.Range("R2").Value = 0.06876867037199999
$ws.Range("S2").Value = 0.08094716512538251
$ws.Range("T2").Value = 0.08094716512538253

# --- Update row 3 (MuSCs -> Fgf8 -> Fgfr2 -> FAPs) with new TPM-derived values ---
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Fgf8"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07200033333333333
$ws.Range("H3").Value = 0.216001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.092289666666667
$ws.Range("N3").Value = 3.276869
$ws.Range("O3").Value = 0.8331551016962769
$ws.Range("P3").Value = 0.833155101696277
$ws.Range("Q3").Value = 0.07864522009655556
$ws.Range("R3").Value = 0.707806980869
$ws.Range("S3").Value = 0.8331551016962769
$ws.Range("T3").Value = 0.833155101696277

# --- Update row 4 (MuSCs -> Fgf8 -> Fgfr2 -> MuSCs) with new TPM-derived values ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Fgf8"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07200033333333333
$ws.Range("H4").Value = 0.216001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1126143333333333
$ws.Range("N4").Value = 0.337843
$ws.Range("O4").Value = 0.08589773317834044
$ws.Range("P4").Value = 0.08589773317834046
$ws.Range("Q4").Value = 0.008108269538111111
$ws.Range("R4").Value = 0.072974425843
$ws.Range("S4").Value = 0.08589773317834044
$ws.Range("T4").Value = 0.08589773317834046

# --- Remove the old rows 5-7 (previously "Inflammatory-Mac" sending-cluster rows) ---
$ws.Rows("5:7").Delete()
